# Apply the two changes described by the diff:
#  1. Slide 7 ("Content Placeholder 2"): strike-through the paragraph
#     "Automate AWS Lambda updates (e.g. data cleaning) using AWS CLI/github
#     actions so that don't have to keep copy and pasting"
#  2. Slide 8 ("TextBox 22"): merge the " role " run into the following
#     "assumed " run, keeping the latter's rPr (dirty="0").

$p = $ppt.ActivePresentation

# --- Change 1: slide 7, paragraph 5 of the content placeholder -> strikethrough
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange
$para5 = $tr7.Paragraphs(5, 1)
$para5.Font.Strikethrough = -1

# --- Change 2: slide 8, TextBox 22 -> merge " role " + "assumed " runs
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(6)
$tr8 = $shp8.TextFrame.TextRange

$fullText = $tr8.Text
$roleIdx = $fullText.IndexOf(" role ")
$roleRange = $tr8.Characters($roleIdx + 1, 6)
$roleRange.Text = ""

$fullText2 = $tr8.Text
$assumedIdx = $fullText2.IndexOf("assumed ")
$assumedRange = $tr8.Characters($assumedIdx + 1, 8)
$assumedRange.Text = " role assumed "
